$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F334").Value = 195128
$ws.Range("G334").Value = 3417
$ws.Range("F335").Value = 130405
$ws.Range("G335").Value = 2974
$ws.Range("F336").Value = 102041
$ws.Range("G336").Value = 3359
$ws.Range("F337").Value = 104278
$ws.Range("F338").Value = 226417
$ws.Range("G338").Value = 3184
$ws.Range("F339").Value = 655755
$ws.Range("G339").Value = 5453
$ws.Range("F340").Value = 380700
$ws.Range("G340").Value = 3260
$ws.Range("F341").Value = 291601
$ws.Range("G341").Value = 3655
$ws.Range("F342").Value = 179101
$ws.Range("G342").Value = 3065
$ws.Range("F343").Value = 132082
$ws.Range("G343").Value = 2961
$ws.Range("F344").Value = 135746
$ws.Range("G344").Value = 2488
$ws.Range("F345").Value = 290969
$ws.Range("G345").Value = 3308
$ws.Range("F346").Value = 668813
$ws.Range("G346").Value = 4779
$ws.Range("F347").Value = 340657
$ws.Range("G347").Value = 2886
$ws.Range("F348").Value = 231754
$ws.Range("G348").Value = 3241
$ws.Range("F349").Value = 159928
$ws.Range("G349").Value = 2750
$ws.Range("F350").Value = 127452
$ws.Range("G350").Value = 2977
$ws.Range("F351").Value = 150341
$ws.Range("G351").Value = 2823
$ws.Range("F352").Value = 306808
$ws.Range("G352").Value = 3545
$ws.Range("F353").Value = 720099
$ws.Range("G353").Value = 5255
$ws.Range("F354").Value = 306803
$ws.Range("G354").Value = 2807
$ws.Range("F355").Value = 222425
$ws.Range("G355").Value = 3446
$ws.Range("F356").Value = 160105
$ws.Range("G356").Value = 2885
$ws.Range("F357").Value = 138168
$ws.Range("G357").Value = 3024
$ws.Range("F358").Value = 157382
$ws.Range("F359").Value = 321030
$ws.Range("G359").Value = 3357
$ws.Range("F360").Value = 741193
$ws.Range("G360").Value = 5085
$ws.Range("F361").Value = 329788
$ws.Range("G361").Value = 2589
$ws.Range("F362").Value = 226035
$ws.Range("G362").Value = 3120
$ws.Range("F363").Value = 186279
$ws.Range("G363").Value = 2743
$ws.Range("F364").Value = 164987
$ws.Range("G364").Value = 2427
$ws.Range("F365").Value = 178545
$ws.Range("G365").Value = 2325
$ws.Range("F366").Value = 332067
$ws.Range("G366").Value = 2798
$ws.Range("F367").Value = 741817
$ws.Range("G367").Value = 3770
$ws.Range("F368").Value = 339366
$ws.Range("G368").Value = 2245
$ws.Range("F369").Value = 227741
$ws.Range("G369").Value = 2508
$ws.Range("F370").Value = 174322
$ws.Range("G370").Value = 1960
$ws.Range("F371").Value = 148124
$ws.Range("G371").Value = 1811
